$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New "K" values (column G) regenerated from s_vals calculation
$kValues = @{
    2 = 0
    3 = 0
    4 = 1
    5 = 1
    6 = 0
    7 = 3
    8 = 0
    9 = 1
    10 = 2
    11 = 1
    12 = 0
    13 = 0
    14 = 2
    15 = 1
    16 = 0
    17 = 1
    18 = 0
    19 = 2
    20 = 1
    21 = 1
    22 = 1
    23 = 0
    24 = 0
    25 = 0
    26 = 1
    27 = 2
    28 = 1
    29 = 1
    30 = 0
    31 = 2
    32 = 1
    33 = 2
    34 = 0
    35 = 3
    36 = 4
    37 = 0
    38 = 1
    39 = 1
    40 = 1
    41 = 0
    42 = 2
    43 = 1
    44 = 1
    45 = 0
    46 = 1
    47 = 1
    48 = 1
    49 = 1
    50 = 2
    51 = 1
    52 = 1
    53 = 4
    54 = 2
    55 = 3
    56 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
